$d = $word.ActiveDocument

# --- 1. Fix in-text citation placeholders -> formatted author-date citations ---
$citations = @(
    @("[@peng2011]", "(Peng, 2011)"),
    @("[@bollen2015]", "(Bollen et al., 2015)"),
    @("[@rosenthal1979]", "(Rosenthal, 1979)"),
    @("[@simmons2011]", "(Simmons et al., 2011)"),
    @("[@young2008]", "(Young et al., 2008)"),
    @("@frisch1933", "Frisch (1933)"),
    @("@mccullough2008", "McCullough et al. (2008)"),
    @("@knuth1992", "Knuth (1992)"),
    @("@ramsey", "Ramsey (n.d.)"),
    @("@gentleman2007", "Gentleman and Lang (2007)")
)

foreach ($pair in $citations) {
    $old = $pair[0]
    $new = $pair[1]
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- 2. Insert the bibliography (References) block ---
# Locate the boundary right after the "Referanser" heading paragraph (where the
# empty "refs" bookmark sits) and before the "Appendiks" heading paragraph.
$appendiksStart = $d.Bookmarks("appendiks").Range.Start
$insertionPoint = $d.Range($appendiksStart, $appendiksStart)

$bibliographyPackageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:bookmarkStart w:id="101" w:name="refs"/><w:bookmarkStart w:id="102" w:name="ref-bollen2015"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Bollen, K., Cacioppo, J. T., Krosnick, J. A., Olds, J. L., and Kaplan, R. M. (2015).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Social, Behavioral, and Economic Sciences Perspectives on Robust and Reliable Science</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">(Report of the Subcommittee on Replicability in Science Advisory Committee to the National Science Foundation Directorate for Social, Behavioral, and Economic Sciences). NSF.</w:t></w:r></w:p><w:bookmarkEnd w:id="102"/><w:bookmarkStart w:id="103" w:name="ref-frisch1933"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Frisch, R. (1933). Editor’s note.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Econometrica</w:t></w:r><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">1</w:t></w:r><w:r><w:t xml:space="preserve">(1), 1–4.</w:t></w:r></w:p><w:bookmarkEnd w:id="103"/><w:bookmarkStart w:id="104" w:name="ref-gentleman2007"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Gentleman, R., and Lang, D. T. (2007). Statistical analyses and reproducible research.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Journal of Computational and Graphical Statistics</w:t></w:r><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">16</w:t></w:r><w:r><w:t xml:space="preserve">(1), 1–23.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId1001"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve">https://doi.org/10.1198/106186007X178663</w:t></w:r></w:hyperlink></w:p><w:bookmarkEnd w:id="104"/><w:bookmarkStart w:id="105" w:name="ref-knuth1992"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Knuth, D. E. (1992).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Literate Programming</w:t></w:r><w:r><w:t xml:space="preserve">. Cambridge University Press.</w:t></w:r></w:p><w:bookmarkEnd w:id="105"/><w:bookmarkStart w:id="106" w:name="ref-mccullough2008"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">McCullough, B. D., McGeary, K. A., and Harrison, T. D. (2008). Do economics journal archives promote replicable research?</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Canadian Journal of Economics/Revue Canadienne d’économique</w:t></w:r><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">41</w:t></w:r><w:r><w:t xml:space="preserve">(4), 1406–1420.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId1002"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve">https://doi.org/10.1111/j.1540-5982.2008.00509.x</w:t></w:r></w:hyperlink></w:p><w:bookmarkEnd w:id="106"/><w:bookmarkStart w:id="107" w:name="ref-peng2011"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Peng, R. D. (2011). Reproducible Research in Computational Science.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Science</w:t></w:r><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">334</w:t></w:r><w:r><w:t xml:space="preserve">(6060), 1226–1227.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId1003"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve">https://doi.org/10.1126/science.1213847</w:t></w:r></w:hyperlink></w:p><w:bookmarkEnd w:id="107"/><w:bookmarkStart w:id="108" w:name="ref-ramsey"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Ramsey, N. (n.d.).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Noweb home page</w:t></w:r><w:r><w:t xml:space="preserve">.</w:t></w:r></w:p><w:bookmarkEnd w:id="108"/><w:bookmarkStart w:id="109" w:name="ref-rosenthal1979"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Rosenthal, R. (1979).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">The file drawer problem and tolerance for null results.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">86</w:t></w:r><w:r><w:t xml:space="preserve">, 638–641.</w:t></w:r></w:p><w:bookmarkEnd w:id="109"/><w:bookmarkStart w:id="110" w:name="ref-simmons2011"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Simmons, J. P., Nelson, L. D., and Simonsohn, U. (2011). False-positive psychology: Undisclosed flexibility in data collection and analysis allows presenting anything as significant.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Psychological Science</w:t></w:r><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">22</w:t></w:r><w:r><w:t xml:space="preserve">(11), 1359–1366.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId1004"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve">https://doi.org/10.1177/0956797611417632</w:t></w:r></w:hyperlink></w:p><w:bookmarkEnd w:id="110"/><w:bookmarkStart w:id="111" w:name="ref-young2008"/><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Young, N. S., Ioannidis, J. P. A., and Al-Ubaydli, O. (2008). Why Current Publication Practices May Distort Science.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">PLoS Medicine</w:t></w:r><w:r><w:t xml:space="preserve">,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">5</w:t></w:r><w:r><w:t xml:space="preserve">(10), e201.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId1005"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve">https://doi.org/10.1371/journal.pmed.0050201</w:t></w:r></w:hyperlink></w:p><w:bookmarkEnd w:id="111"/><w:bookmarkEnd w:id="101"/></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1001" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1198/106186007X178663" TargetMode="External"/><Relationship Id="rId1002" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1111/j.1540-5982.2008.00509.x" TargetMode="External"/><Relationship Id="rId1003" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1126/science.1213847" TargetMode="External"/><Relationship Id="rId1004" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1177/0956797611417632" TargetMode="External"/><Relationship Id="rId1005" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1371/journal.pmed.0050201" TargetMode="External"/></Relationships></pkg:xmlData></pkg:part></pkg:package>'

$null = $insertionPoint.InsertXML($bibliographyPackageXml)

Write-Output "edit applied"
